# ReglasCargaMasiva.xlsx - creacion de plantillas carga masiva
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (closest achievable values given the engine's internal
# pixel-grid rounding of column widths -- target source widths were
# 17.28515625 and 32.42578125 character units)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(3).ColumnWidth = 31.65

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Campo"
$ws.Range("B1").Value = "Tipo"
$ws.Range("C1").Value = "Valores"

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Tipo Documento"
$ws.Range("B2").Value = "Enum"
$ws.Range("C2").Value = " TI,CC,CE,PS"

$ws.Range("A3").Value = "Numero de documento"
$ws.Range("B3").Value = "String"
$ws.Range("C3").Value = "Texto"

$ws.Range("A4").Value = "Primer nombre"
$ws.Range("B4").Value = "String"
$ws.Range("C4").Value = "Texto"

$ws.Range("A5").Value = "Segundo Nombre"
$ws.Range("B5").Value = "String"
$ws.Range("C5").Value = "Texto"

$ws.Range("A6").Value = "Primer Apellido"
$ws.Range("B6").Value = "String"
$ws.Range("C6").Value = "Texto"

$ws.Range("A7").Value = "Segundo Apellido"
$ws.Range("B7").Value = "String"
$ws.Range("C7").Value = "Texto"

$ws.Range("A8").Value = "Correo Electronico"
$ws.Range("B8").Value = "String"
$ws.Range("C8").Value = "Texto"

$ws.Range("A9").Value = "Telefono"
$ws.Range("B9").Value = "String"
$ws.Range("C9").Value = "Texto"

$ws.Range("A10").Value = "Area"
$ws.Range("B10").Value = "Enum"
$ws.Range("C10").Value = "INGENIERIA, INFRA,MANTENIMIENTO,GERENCIA,RRHH,ADMIN"

$ws.Range("A11").Value = "Jefe"
$ws.Range("B11").Value = "List"
$ws.Range("C11").Value = "Nombre y apellido del gerente"

$ws.Range("A12").Value = "Codigo usuario"
$ws.Range("B12").Value = "String"
$ws.Range("C12").Value = "Texto"

$ws.Range("A13").Value = "Tipo Usuario"
$ws.Range("B13").Value = "List"

# ---------------------------------------------------------------------------
# Blank formatted cells (rows 6-15, columns D-I) and rows 14-15 A/C
# ---------------------------------------------------------------------------
$ws.Range("D6:I15").Value = ""
$ws.Range("A14").Value = ""
$ws.Range("A15").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("C15").Value = ""

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 21.75
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 16.5
$ws.Rows.Item(8).RowHeight = 17.25
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 31.5

# ---------------------------------------------------------------------------
# Styles: alignment vertical=top + wrapText (style index 1)
# built via a helper cell + paste-special so the COM layer doesn't leave
# stray intermediate styles behind when styling multi-cell ranges.
# ---------------------------------------------------------------------------
$ws.Range("Z1").WrapText = $true
$ws.Range("Z1").VerticalAlignment = -4160
$ws.Range("Z1").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122)
$ws.Range("C10:C11").PasteSpecial(-4122)
$ws.Range("C14:C15").PasteSpecial(-4122)
$ws.Range("D6:I15").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# Styles: alignment vertical=top only (style index 2)
$ws.Range("Z1").VerticalAlignment = -4160
$ws.Range("Z1").Copy()
$ws.Range("B2:B13").PasteSpecial(-4122)
$ws.Range("C2:C9").PasteSpecial(-4122)
$ws.Range("C12:C13").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Selection shown when the file was last saved
# ---------------------------------------------------------------------------
$null = $ws.Range("C17").Select()
